$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark "Shortest Route with Line Change Penalties (Dijkstras)" (row 9) and its
# sub-items (rows 10-12) as Attempted (checkbox / boolean column C).
$ws.Range("C9").Value = $true
$ws.Range("C10").Value = $true
$ws.Range("C11").Value = $true
$ws.Range("C12").Value = $true

# Update the active selection to reflect where the user last clicked.
$ws.Range("C9").Select()

$wb.Application.Calculate()
